$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$np = $s.NotesPage
$ph = $np.Shapes.AddPlaceholder(2)
$ph.TextFrame.TextRange.Text = "Speaker notes here"
